$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "box"
$ws.Range("C2").Value = 30
$ws.Range("D2").Value = 0.89

$ws.Range("B3").Value = "core"
$ws.Range("C3").Value = 1500
$ws.Range("D3").Value = 0.023

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "label"
$ws.Range("C4").Value = 60
$ws.Range("D4").Value = 799
